$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.424.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7050"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07848"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.38%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3133"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08014"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.890.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.193"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.7002"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.439"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.503.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008323"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.142.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.601"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.86%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1555"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.016"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.319"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.268"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.203"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05297"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.882"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7497"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.164"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01872"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.259.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.740"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8988"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.65"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.957"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.35%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("E46").Value = "  -0.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.038.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5184"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.787"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.488"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4303"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.76%  "
